$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.160006403923035
$ws.Range("B1").Value = 2.407143831253052
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.384583950042725
$ws.Range("E1").Value = 1.229294657707214
